$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 170; this shifts the existing rows
# 170-199 down to 172-201 (formatting/style carries with the shift).
$ws.Rows("170:171").Insert()

# --- New row 170: Brócoli, Primera ---
$ws.Cells.Item(170, 1).Value = 1
$ws.Cells.Item(170, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(170, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(170, 4).Value = 44474
$ws.Cells.Item(170, 5).Value = 15
$ws.Cells.Item(170, 6).Value = 100112023
$ws.Cells.Item(170, 7).Value = "Brócoli"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 1000
$ws.Cells.Item(170, 11).Value = 600
$ws.Cells.Item(170, 12).Value = 700
$ws.Cells.Item(170, 13).Value = 650
$ws.Cells.Item(170, 14).Value = "`$/unidad"
$ws.Cells.Item(170, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(170, 16).Value = 650
$ws.Cells.Item(170, 17).Value = 1
$ws.Cells.Item(170, 18).Value = "Hortaliza"

# --- New row 171: Brócoli, Segunda ---
$ws.Cells.Item(171, 1).Value = 1
$ws.Cells.Item(171, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(171, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(171, 4).Value = 44474
$ws.Cells.Item(171, 5).Value = 15
$ws.Cells.Item(171, 6).Value = 100112023
$ws.Cells.Item(171, 7).Value = "Brócoli"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Segunda"
$ws.Cells.Item(171, 10).Value = 1300
$ws.Cells.Item(171, 11).Value = 400
$ws.Cells.Item(171, 12).Value = 500
$ws.Cells.Item(171, 13).Value = 450
$ws.Cells.Item(171, 14).Value = "`$/unidad"
$ws.Cells.Item(171, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(171, 16).Value = 450
$ws.Cells.Item(171, 17).Value = 1
$ws.Cells.Item(171, 18).Value = "Hortaliza"

Write-Host "Dimension after edit:" $ws.UsedRange.Address()
